$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of collection data (MCH147) ---
$ws.Cells.Item(2,1).Value = "MCH147"
$ws.Cells.Item(2,3).Value = "SOLOMON MAHLANGU HERO OF THE FREEDOM STRUGGLE, SOMAFCO"
$ws.Cells.Item(2,5).Value = "Series"
$ws.Cells.Item(2,6).Value = "1 Box"
$ws.Cells.Item(2,7).Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# Apply the same font (Calibri 10pt, text1 theme color) used for the
# rest of the data row to every cell that belongs to the new record,
# including the two blank cells (date_s / file_path) that still carry
# the row's formatting even though they hold no value.
foreach ($col in 1,3,4,5,6,7,8) {
    $cell = $ws.Cells.Item(2, $col)
    $cell.Font.Name = "Calibri"
    $cell.Font.ThemeColor = 1
}

# Match the row heights already used on the header row
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75

# Re-apply the frozen header row / current selection now that the sheet
# has a second row of real data
[void]$ws.Range("A2:J2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
